$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# --- Header row (B1:E1) ---
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# --- Row 2 data (B2:E2) ---
$ws.Range("B2").Value = 40.018403858124941
$ws.Range("C2").Value = 24.832978138124986
$ws.Range("D2").Value = 35.186674998124943
$ws.Range("E2").Value = 31.556935828125006

# --- Row 3 data (B3:E3) ---
$ws.Range("B3").Value = 34.57344337125005
$ws.Range("C3").Value = 20.167593688124953
$ws.Range("D3").Value = 35.947099781250074
$ws.Range("E3").Value = 35.443980539999927

# --- Selection range update ---
$ws.Range("B1:E3").Select()
